$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for Iteration 5 and 6, copying style from existing header cell (F1)
$ws.Range("G1").Value = "Iteration - 5"
$ws.Range("H1").Value = "Iteration - 6"
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Row 2 (Mean) - update existing values and add new ones
$ws.Range("B2").Value = 0.2225284506110743
$ws.Range("C2").Value = -0.01523612079177302
$ws.Range("D2").Value = 0.03247825077474167
$ws.Range("E2").Value = 0.03859681228091425
$ws.Range("F2").Value = 0.004026113512708367
$ws.Range("G2").Value = 0.004998000895654252
$ws.Range("H2").Value = 0.006053141156931599

# Row 3 (Standard Deviation) - update existing values and add new ones
$ws.Range("B3").Value = 1.666541573859614
$ws.Range("C3").Value = 0.4868097093947357
$ws.Range("D3").Value = 0.2279358809685001
$ws.Range("E3").Value = 0.1898772759621833
$ws.Range("F3").Value = 0.1235793968233901
$ws.Range("G3").Value = 0.1098257311943792
$ws.Range("H3").Value = 0.1004375325099085

# Row 4 (Outlier) - update existing values and add new ones
$ws.Range("B4").Value = 386
$ws.Range("C4").Value = 107
$ws.Range("D4").Value = 122
$ws.Range("E4").Value = 440
$ws.Range("F4").Value = 122
$ws.Range("G4").Value = 152
$ws.Range("H4").Value = 64
